$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.302.82"
$ws.Range("E2").Value = "  -7.20%  "

# Row 3
$ws.Range("D3").Value = "3.543.08"
$ws.Range("E3").Value = "  -3.95%  "

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "391.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.02%  "

# Row 7
$ws.Range("D7").Value = "3.534.88"
$ws.Range("E7").Value = "  -3.92%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.584"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -10.63%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.676"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -12.02%  "

# Row 11
$ws.Range("E11").Value = "  -22.92%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000324"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -26.64%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.03%  "

# Row 14
$ws.Range("D14").Value = "4.102.50"
$ws.Range("E14").Value = "  -3.58%  "

# Row 15
$ws.Range("E15").Value = "  -7.26%  "

# Row 16
$ws.Range("E16").Value = "  -2.93%  "

# Row 17
$ws.Range("D17").Value = "3.558.55"
$ws.Range("E17").Value = "  -2.78%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.71%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.31%  "

# Row 20
$ws.Range("D20").Value = "63.271.31"
$ws.Range("E20").Value = "  -7.17%  "

# Row 21
$ws.Range("E21").Value = "  -10.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -14.72%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.77%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.62%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.44"
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.73%  "

# Row 28
$ws.Range("E28").Value = "  -8.65%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -13.90%  "

# Row 30
$ws.Range("E30").Value = "  -5.34%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.02%  "

# Row 32
$ws.Range("E32").Value = "  -7.03%  "

# Row 33
$ws.Range("E33").Value = "  -5.47%  "

# Row 34
$ws.Range("E34").Value = "  -5.98%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.12%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "36.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.25%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.67%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0435"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -11.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.997"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.13%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0648"
$ws.Range("E40").Value = "  -17.23%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.97%  "

# Row 42
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.130"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -13.20%  "

# Row 43
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.64%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "140.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.01%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +16.34%  "

# Row 46
$ws.Range("E46").Value = "  -1.14%  "

# Row 47
$ws.Range("E47").Value = "  -7.13%  "

# Row 48
$ws.Range("E48").Value = "  -4.63%  "

# Row 49
$ws.Range("E49").Value = "  -8.78%  "

# Row 50
$ws.Range("E50").Value = "  -9.83%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.274"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.68%  "
